$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 41.1
$ws.Range("N3").Value = 24.66
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 3
$ws.Range("U3").Value = 3

# Row 4
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("P4").Value = 14
$ws.Range("Q4").Value = 0
$ws.Range("U4").Value = 0
$ws.Rows.Item(4).EntireRow.Hidden = $true

# Row 5
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("P5").Value = 6
$ws.Range("Q5").Value = 0
$ws.Range("U5").Value = 0
$ws.Rows.Item(5).EntireRow.Hidden = $true

# Row 6
$ws.Range("M6").Value = 7.9
$ws.Range("N6").Value = 4.74
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 4

# Row 7
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = 1.2
$ws.Range("P7").Value = 7
$ws.Range("Q7").Value = 1
$ws.Range("U7").Value = 1

# Row 8
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("P8").Value = 7
$ws.Range("Q8").Value = 0
$ws.Range("U8").Value = 0
$ws.Rows.Item(8).EntireRow.Hidden = $true

# Row 9
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("P9").Value = 12
$ws.Range("Q9").Value = 0
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 2
$ws.Range("U9").Value = 0
$ws.Rows.Item(9).EntireRow.Hidden = $true

# Row 10
$ws.Range("M10").Value = 9.88
$ws.Range("N10").Value = 5.93
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 2
$ws.Range("U10").Value = 5

# Row 11
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("P11").Value = 8
$ws.Range("Q11").Value = 0
$ws.Range("S11").Value = 2
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 0
$ws.Rows.Item(11).EntireRow.Hidden = $true

# Row 12
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("P12").Value = 9
$ws.Range("Q12").Value = 0
$ws.Range("U12").Value = 0
$ws.Rows.Item(12).EntireRow.Hidden = $true

# Row 13
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("P13").Value = 3
$ws.Range("Q13").Value = 0
$ws.Range("U13").Value = 0
$ws.Rows.Item(13).EntireRow.Hidden = $true

# Row 14
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("P14").Value = 5
$ws.Range("Q14").Value = 0
$ws.Range("U14").Value = 0
$ws.Rows.Item(14).EntireRow.Hidden = $true

# Row 15
$ws.Range("M15").Value = 8
$ws.Range("N15").Value = 4.8
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 4

# Row 16
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("P16").Value = 9
$ws.Range("Q16").Value = 0
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 0
$ws.Rows.Item(16).EntireRow.Hidden = $true

# Row 19
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("P19").Value = 9
$ws.Range("Q19").Value = 0
$ws.Range("U19").Value = 0
$ws.Rows.Item(19).EntireRow.Hidden = $true

# Row 20
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("P20").Value = 6
$ws.Range("Q20").Value = 0
$ws.Range("U20").Value = 0
$ws.Rows.Item(20).EntireRow.Hidden = $true

# Row 21
$ws.Range("M21").Value = 7.9
$ws.Range("N21").Value = 4.74
$ws.Range("S21").Value = 1
$ws.Range("T21").Value = 1
$ws.Range("U21").Value = 4

# Row 22
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("P22").Value = 17
$ws.Range("Q22").Value = 0
$ws.Range("U22").Value = 0
$ws.Rows.Item(22).EntireRow.Hidden = $true

# Row 23
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("P23").Value = 2
$ws.Range("Q23").Value = 0
$ws.Range("U23").Value = 0
$ws.Rows.Item(23).EntireRow.Hidden = $true

# Row 24
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("P24").Value = 12
$ws.Range("Q24").Value = 0
$ws.Range("U24").Value = 0
$ws.Rows.Item(24).EntireRow.Hidden = $true

# Row 25
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("P25").Value = 3
$ws.Range("Q25").Value = 0
$ws.Range("U25").Value = 0
$ws.Rows.Item(25).EntireRow.Hidden = $true

# Row 26
$ws.Range("M26").Value = 4
$ws.Range("N26").Value = 2.4
$ws.Range("P26").Value = 2
$ws.Range("Q26").Value = 1
$ws.Range("S26").Value = 1
$ws.Range("T26").Value = 1
$ws.Range("U26").Value = 2

# Row 27
$ws.Range("M27").Value = 9.88
$ws.Range("N27").Value = 5.93
$ws.Range("S27").Value = 3
$ws.Range("T27").Value = 3
$ws.Range("U27").Value = 5

# Row 29
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("P29").Value = 8
$ws.Range("Q29").Value = 0
$ws.Range("S29").Value = 1
$ws.Range("T29").Value = 1
$ws.Range("U29").Value = 0
$ws.Rows.Item(29).EntireRow.Hidden = $true

# Row 30
$ws.Range("L30").Value = 0

# Row 31
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("P31").Value = 10
$ws.Range("Q31").Value = 0
$ws.Range("U31").Value = 0
$ws.Rows.Item(31).EntireRow.Hidden = $true

# Row 32
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("P32").Value = 30
$ws.Range("Q32").Value = 0
$ws.Range("U32").Value = 0
$ws.Rows.Item(32).EntireRow.Hidden = $true

# Row 33
$ws.Range("M33").Value = 14
$ws.Range("N33").Value = 8.4
$ws.Range("P33").Value = 1
$ws.Range("Q33").Value = 7
$ws.Range("U33").Value = 7

# Row 34
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("P34").Value = 3
$ws.Range("Q34").Value = 0
$ws.Range("S34").Value = 1
$ws.Range("T34").Value = 1
$ws.Range("U34").Value = 0
$ws.Rows.Item(34).EntireRow.Hidden = $true

# Row 35
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("P35").Value = 7
$ws.Range("Q35").Value = 0
$ws.Range("S35").Value = 1
$ws.Range("T35").Value = 1
$ws.Range("U35").Value = 0
$ws.Rows.Item(35).EntireRow.Hidden = $true

# Row 36
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("P36").Value = 19
$ws.Range("Q36").Value = 0
$ws.Range("S36").Value = 1
$ws.Range("T36").Value = 1
$ws.Range("U36").Value = 0
$ws.Rows.Item(36).EntireRow.Hidden = $true

# Row 37
$ws.Range("P37").Value = 1
$ws.Range("Q37").Value = 7
$ws.Range("S37").Value = 1
$ws.Range("T37").Value = 1

# Row 38
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("P38").Value = 12
$ws.Range("Q38").Value = 0
$ws.Range("U38").Value = 0
$ws.Rows.Item(38).EntireRow.Hidden = $true

# Row 39
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("P39").Value = 5
$ws.Range("Q39").Value = 0
$ws.Range("U39").Value = 0
$ws.Rows.Item(39).EntireRow.Hidden = $true

# Summary metrics (rows 42, 44, 53)
$ws.Range("C42").Value = 55
# Leading apostrophe forces text entry so '144.51€' is not auto-coerced to a currency number
$ws.Range("C44").Value = "'144.51€"
$ws.Range("C53").Value = -1
